$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the old "UART 1" row (row 11) for the new
# EXTI / PB1 (DIO2) and EXTI / PB2 (DIO0) pin mappings.
$ws.Rows("11:12").Insert()

# Fill column A first (row labels), then column B, then C, then D -- this
# mirrors how the table was authored and keeps the shared-string order
# consistent with the source workbook.
$ws.Range("A11").Value = "EXTI / PB1"
$ws.Range("A12").Value = "EXTI / PB2"

$ws.Range("B11").Value = "DIO2"
$ws.Range("B12").Value = "DIO0"

$ws.Range("C11").Value = "-"
$ws.Range("C12").Value = "-"

$ws.Range("D11").Value = "-"
$ws.Range("D12").Value = "-"

# The GPS VIN row no longer marks the 4th (UART in) column with 5V.
$ws.Range("D4").Value = "-"

# Append the new RESET / 1PPS pin rows at the bottom of the table.
$ws.Range("A19").Value = "PB4"
$ws.Range("A20").Value = "PB5"

$ws.Range("C20").Value = "RESET"
$ws.Range("C19").Value = "1PPS"

$ws.Range("C19").HorizontalAlignment = -4108
$ws.Range("C20").HorizontalAlignment = -4108

$ws.Range("C19").Select()
